$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge T1:T2 before clearing (it was previously merged)
$ws.Range("T1:T2").UnMerge()

# Clear T column except row 6 (new formulas go there), and clear old T103 totals row
$ws.Range("T1:T5").Clear()
$ws.Range("T7:T103").Clear()
$ws.Range("A103:Z103").Clear()

# Add new S6/T6 formulas
$ws.Range("S6").Formula = "=0.9*((R6-D6)^0.51)*((SQRT((P6-B6)^2+(Q6-C6)^2)^(-0.35)))"
$ws.Range("T6").Formula = "=IF(M6=0,TRUE,OR(AND(S6<(M6*1.03),S6>(M6*0.97)),O6))"

# --- Conditional formatting updates ---

# Re-create K3:K102 rules (causes fresh dxf entries, matching diff's re-indexed dxfId 8/7)
$fcsK = $ws.Range("K3:K102").FormatConditions
$fcsK.Delete()
$rKTrue = $fcsK.Add(1, 3, "=TRUE")
$rKTrue.Font.Color = 24832
$rKTrue.Interior.Color = 13561798
$rKFalse = $fcsK.Add(1, 3, "=FALSE")
$rKFalse.Font.Color = 393372
$rKFalse.Interior.Color = 13551615

# Add new FALSE rule to O3:O45 (reuses existing red dxf at index 1)
$fcsO = $ws.Range("O3:O45").FormatConditions
$rOFalse = $fcsO.Add(1, 3, "=FALSE")
$rOFalse.Font.Color = 393372
$rOFalse.Interior.Color = 13551615

# Add new TRUE rule to T6 (green)
$fcsT6 = $ws.Range("T6").FormatConditions
$rT6True = $fcsT6.Add(1, 3, "=TRUE")
$rT6True.Font.Color = 24832
$rT6True.Interior.Color = 13561798

# Set explicit priorities to match target ordering
$rKTrue.Priority = 4
$rKFalse.Priority = 5
$fcsO.Item(1).Priority = 3
$rOFalse.Priority = 1
$rT6True.Priority = 2

Write-Host "Done"
